$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K6").Value = 4291.1538
$ws.Range("I6").Value = 1430.3846
$ws.Range("M6").Value = -4179.1538
$ws.Range("H6").Value = 1430.3846
$ws.Range("K86").Value = 333340000
$ws.Range("J86").Value = 222227230
$ws.Range("I86").Value = 333340000
$ws.Range("M86").Value = -333338877
$ws.Range("H86").Value = 277783600
$ws.Range("N86").Value = -222229476
$ws.Range("L86").Value = 222227230
$ws.Range("J89").Value = 222227230
$ws.Range("K89").Value = 1666700000
$ws.Range("I89").Value = 333340000
$ws.Range("M89").Value = -1666694384
$ws.Range("N89").Value = -1111147382
$ws.Range("H89").Value = 277783600
$ws.Range("L89").Value = 1111136150
$ws.Range("J125").Value = 2025
$ws.Range("K125").Value = 20470752
$ws.Range("I125").Value = 2274528
$ws.Range("M125").Value = -20468292
$ws.Range("N125").Value = -23145
$ws.Range("H125").Value = 1895777.5
$ws.Range("L125").Value = 18225
$ws.Range("I137").Value = 8494.632
$ws.Range("M137").Value = -22933.896
$ws.Range("H137").Value = 3581549.2
$ws.Range("K137").Value = 25483.896

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1381
$ws.Range("K2").Value = 1381
$ws.Range("M2").Value = -1268
$ws.Range("H2").Value = 1381
$ws.Range("H34").Value = 254750
$ws.Range("L34").Value = 500000
$ws.Range("J34").Value = 500000
$ws.Range("N34").Value = -500542
$ws.Range("K74").Value = 1649.925
$ws.Range("J74").Value = 1484770.6
$ws.Range("I74").Value = 1649.925
$ws.Range("M74").Value = -775.925
$ws.Range("N74").Value = -1486518.6
$ws.Range("H74").Value = 406137.4
$ws.Range("L74").Value = 1484770.6
$ws.Range("H77").Value = 406137.4
$ws.Range("L77").Value = 7423853
$ws.Range("J77").Value = 1484770.6
$ws.Range("I77").Value = 1649.925
$ws.Range("K77").Value = 8249.625
$ws.Range("M77").Value = -3881.625
$ws.Range("N77").Value = -7432589
$ws.Range("K88").Value = 1253
$ws.Range("I88").Value = 1253
$ws.Range("M88").Value = -847
$ws.Range("H88").Value = 2176.5386
$ws.Range("M91").Value = 151
$ws.Range("H91").Value = 2176.5386
$ws.Range("K91").Value = 1253
$ws.Range("I91").Value = 1253
$ws.Range("I102").Value = 1986.375
$ws.Range("K102").Value = 1986.375
$ws.Range("M102").Value = -364.375
$ws.Range("H102").Value = 2182.4614
$ws.Range("K116").Value = 1381
$ws.Range("I116").Value = 1381
$ws.Range("M116").Value = 913
$ws.Range("H116").Value = 1381
$ws.Range("N132").Value = -15129.2501
$ws.Range("H132").Value = 2174.25
$ws.Range("L132").Value = 10069.2501
$ws.Range("J132").Value = 3356.4167

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 1381
$ws.Range("M3").Value = -1267
$ws.Range("H3").Value = 1381
$ws.Range("K3").Value = 1381
$ws.Range("J80").Value = 798.6
$ws.Range("N80").Value = -2794.6
$ws.Range("H80").Value = 83334200
$ws.Range("L80").Value = 798.6
$ws.Range("J83").Value = 798.6
$ws.Range("N83").Value = -13977
$ws.Range("H83").Value = 83334200
$ws.Range("L83").Value = 3993
$ws.Range("K86").Value = 3284.4119
$ws.Range("J86").Value = 4764.5264
$ws.Range("I86").Value = 3284.4119
$ws.Range("M86").Value = -2161.4119
$ws.Range("H86").Value = 4065.5833
$ws.Range("N86").Value = -7010.5264
$ws.Range("L86").Value = 4764.5264
$ws.Range("J89").Value = 4764.5264
$ws.Range("K89").Value = 16422.0595
$ws.Range("I89").Value = 3284.4119
$ws.Range("M89").Value = -10806.0595
$ws.Range("N89").Value = -35054.632
$ws.Range("H89").Value = 4065.5833
$ws.Range("L89").Value = 23822.632
$ws.Range("I105").Value = 11798.546
$ws.Range("K105").Value = 11798.546
$ws.Range("M105").Value = -10051.546
$ws.Range("H105").Value = 10463.305
$ws.Range("M134").Value = -2678.0358
$ws.Range("H134").Value = 23078628
$ws.Range("L134").Value = 245459400
$ws.Range("J134").Value = 81819800
$ws.Range("K134").Value = 5213.0358
$ws.Range("I134").Value = 1737.6786
$ws.Range("N134").Value = -245464470

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 20409424
$ws.Range("M16").Value = -20409137
$ws.Range("H16").Value = 11912847
$ws.Range("K16").Value = 20409424
$ws.Range("I105").Value = 2384.25
$ws.Range("K105").Value = 2384.25
$ws.Range("M105").Value = -637.25
$ws.Range("H105").Value = 2384.25
$ws.Range("M113").Value = -20407254
$ws.Range("H113").Value = 11912847
$ws.Range("K113").Value = 20409424
$ws.Range("I113").Value = 20409424
$ws.Range("M134").Value = -2321.6667
$ws.Range("H134").Value = 2183.125
$ws.Range("K134").Value = 4856.6667
$ws.Range("I134").Value = 1618.8889

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N113").Value = -7946
$ws.Range("H113").Value = 1020.8
$ws.Range("L113").Value = 3606
$ws.Range("J113").Value = 1202
$ws.Range("J121").Value = 17981.182
$ws.Range("N121").Value = -56563.546
$ws.Range("H121").Value = 14493.929
$ws.Range("L121").Value = 53943.546
$ws.Range("I136").Value = 3968.6
$ws.Range("M136").Value = -6805.799999999999
$ws.Range("N136").Value = -52798.5
$ws.Range("H136").Value = 9084.049999999999
$ws.Range("L136").Value = 42598.5
$ws.Range("J136").Value = 14199.5
$ws.Range("K136").Value = 11905.8
$ws.Range("N138").Value = -20030
$ws.Range("H138").Value = 3091.318
$ws.Range("L138").Value = 9750
$ws.Range("J138").Value = 3250
$ws.Range("K139").Value = 22064460
$ws.Range("I139").Value = 7354820
$ws.Range("M139").Value = -22059320
$ws.Range("H139").Value = 6252847

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M7").Value = -499888
$ws.Range("N7").Value = -2000224
$ws.Range("H7").Value = 1250000
$ws.Range("L7").Value = 2000000
$ws.Range("J7").Value = 2000000
$ws.Range("K7").Value = 500000
$ws.Range("I7").Value = 500000
$ws.Range("H8").Value = 1250000
$ws.Range("L8").Value = 2000000
$ws.Range("K8").Value = 500000
$ws.Range("J8").Value = 2000000
$ws.Range("I8").Value = 500000
$ws.Range("M8").Value = -499861
$ws.Range("N8").Value = -2000278
$ws.Range("J11").Value = 3022502
$ws.Range("K11").Value = 20998.625
$ws.Range("I11").Value = 20998.625
$ws.Range("M11").Value = -20859.625
$ws.Range("N11").Value = -3022780
$ws.Range("H11").Value = 1021499.75
$ws.Range("L11").Value = 3022502
$ws.Range("J80").Value = 71484630
$ws.Range("K80").Value = 113882.05
$ws.Range("I80").Value = 113882.05
$ws.Range("M80").Value = -112884.05
$ws.Range("N80").Value = -71486626
$ws.Range("H80").Value = 9909867
$ws.Range("L80").Value = 71484630
$ws.Range("J83").Value = 71484630
$ws.Range("K83").Value = 569410.25
$ws.Range("I83").Value = 113882.05
$ws.Range("M83").Value = -564418.25
$ws.Range("N83").Value = -357433134
$ws.Range("H83").Value = 9909867
$ws.Range("L83").Value = 357423150
$ws.Range("J101").Value = 29666.666
$ws.Range("N101").Value = -36156.666
$ws.Range("H101").Value = 29666.666
$ws.Range("L101").Value = 29666.666
$ws.Range("J107").Value = 1352.8
$ws.Range("N107").Value = -5192.8
$ws.Range("H107").Value = 91748.63
$ws.Range("L107").Value = 1352.8
$ws.Range("N117").Value = -66884
$ws.Range("H117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("I126").Value = 6325
$ws.Range("M126").Value = -16505
$ws.Range("H126").Value = 6325
$ws.Range("L126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 18975
$ws.Range("N126").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I82").Value = 924.5
$ws.Range("K82").Value = 924.5
$ws.Range("M82").Value = -563.5
$ws.Range("H82").Value = 921.7778
$ws.Range("H85").Value = 921.7778
$ws.Range("K85").Value = 924.5
$ws.Range("I85").Value = 924.5
$ws.Range("M85").Value = 323.5
$ws.Range("K100").Value = 2569.5386
$ws.Range("I100").Value = 2569.5386
$ws.Range("M100").Value = -2028.5386
$ws.Range("H100").Value = 2826.9333
$ws.Range("J101").Value = 22615.143
$ws.Range("N101").Value = -29105.143
$ws.Range("H101").Value = 22615.143
$ws.Range("L101").Value = 22615.143
$ws.Range("J122").Value = 4850
$ws.Range("K122").Value = 8578.0905
$ws.Range("I122").Value = 2859.3635
$ws.Range("M122").Value = -6128.0905
$ws.Range("N122").Value = -19450
$ws.Range("H122").Value = 3390.2
$ws.Range("L122").Value = 14550
$ws.Range("N131").Value = -60080
$ws.Range("H131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("J131").Value = 50000

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J2").Value = 3343335.2
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("N2").Value = -3343559.2
$ws.Range("H2").Value = 3343335.2
$ws.Range("L2").Value = 3343335.2
$ws.Range("J107").Value = 1787874.4
$ws.Range("I107").Value = 1957.3529
$ws.Range("K107").Value = 5872.0587
$ws.Range("M107").Value = -3952.0587
$ws.Range("N107").Value = -5367463.199999999
$ws.Range("H107").Value = 867856.5
$ws.Range("L107").Value = 5363623.199999999
$ws.Range("M2").ClearContents()
